$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-18"

# Update the header label shared string text (cell I1 uses this string)
$ws.Range("I1").Value = "2022 (through 10-18)"

# Update the monthly data values for October (row 10) and November (row 11)
$ws.Range("I10").Value = 145
$ws.Range("I11").Value = 60

# Update the Total row to reflect the new sum
$ws.Range("I14").Value = 1337
